$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 2).Value = 0
$ws.Cells.Item(1, 3).Value = 1
$ws.Cells.Item(1, 4).Value = 2
$ws.Cells.Item(1, 5).Value = 3
$ws.Cells.Item(1, 6).Value = 4
$ws.Cells.Item(1, 7).Value = 5
$ws.Cells.Item(1, 8).Value = 6
$ws.Cells.Item(1, 9).Value = 7
$ws.Cells.Item(1, 10).Value = 8
$ws.Cells.Item(1, 11).Value = 9
$ws.Cells.Item(1, 12).Value = 10
$ws.Cells.Item(1, 13).Value = 11
$ws.Cells.Item(1, 14).Value = 12
$ws.Cells.Item(1, 15).Value = 13
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Cells.Item(1, 18).Value = 16
$ws.Cells.Item(1, 19).Value = 17
$ws.Cells.Item(1, 20).Value = 18
$ws.Cells.Item(1, 21).Value = 19
$ws.Cells.Item(1, 22).Value = 20
$ws.Cells.Item(1, 23).Value = 21
$ws.Cells.Item(1, 24).Value = 22
$ws.Cells.Item(1, 25).Value = 23
$ws.Cells.Item(1, 26).Value = 24
$ws.Cells.Item(1, 27).Value = 25
$ws.Cells.Item(1, 28).Value = 26
$ws.Cells.Item(1, 29).Value = 27
$ws.Cells.Item(1, 30).Value = 28
$ws.Cells.Item(1, 31).Value = 29

# Copy the header-row cell format (bold, border, centered) from X1 onto
# the newly added header cells Y1:AE1 so they match the rest of row 1.
$ws.Range("X1").Copy()
$ws.Range("Y1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Cells.Item(2, 2).Value = 0.6551724137931034
$ws.Cells.Item(2, 3).Value = 0.4795918367346939
$ws.Cells.Item(2, 4).Value = 0.5567010309278351
$ws.Cells.Item(2, 5).Value = 0.6703296703296704
$ws.Cells.Item(2, 6).Value = 0.6994818652849741
$ws.Cells.Item(2, 7).Value = 0.7411764705882353
$ws.Cells.Item(2, 8).Value = 0.5652173913043478
$ws.Cells.Item(2, 9).Value = 0.4693877551020407
$ws.Cells.Item(2, 10).Value = 0.5789473684210527
$ws.Cells.Item(2, 11).Value = 0.5402298850574713
$ws.Cells.Item(2, 12).Value = 0.5491329479768786
$ws.Cells.Item(2, 13).Value = 0.4470046082949309
$ws.Cells.Item(2, 14).Value = 0.4619883040935673
$ws.Cells.Item(2, 15).Value = 0.468599033816425
$ws.Cells.Item(2, 16).Value = 0.4380952380952381
$ws.Cells.Item(2, 17).Value = 0.53
$ws.Cells.Item(2, 18).Value = 0.4672131147540983
$ws.Cells.Item(2, 19).Value = 0.5124378109452737
$ws.Cells.Item(2, 20).Value = 0.4764397905759162
$ws.Cells.Item(2, 21).Value = 0.6157635467980295
$ws.Cells.Item(2, 22).Value = 0.6972972972972973
$ws.Cells.Item(2, 23).Value = 0.663551401869159
$ws.Cells.Item(2, 24).Value = 0.6393442622950819
$ws.Cells.Item(2, 25).Value = 0.5436893203883496
$ws.Cells.Item(2, 26).Value = 0.592039800995025
$ws.Cells.Item(2, 27).Value = 0.4975369458128078
$ws.Cells.Item(2, 28).Value = 0.5376344086021506
$ws.Cells.Item(2, 29).Value = 0.5752212389380531
$ws.Cells.Item(2, 30).Value = 0.5692307692307692
$ws.Cells.Item(2, 31).Value = 0.4956521739130436
# Row 3
$ws.Cells.Item(3, 2).Value = 0.3448275862068965
$ws.Cells.Item(3, 3).Value = 0.5204081632653061
$ws.Cells.Item(3, 4).Value = 0.4432989690721649
$ws.Cells.Item(3, 5).Value = 0.3296703296703297
$ws.Cells.Item(3, 6).Value = 0.3005181347150259
$ws.Cells.Item(3, 7).Value = 0.2588235294117647
$ws.Cells.Item(3, 8).Value = 0.4347826086956522
$ws.Cells.Item(3, 9).Value = 0.5306122448979592
$ws.Cells.Item(3, 10).Value = 0.4210526315789473
$ws.Cells.Item(3, 11).Value = 0.4597701149425288
$ws.Cells.Item(3, 12).Value = 0.4508670520231214
$ws.Cells.Item(3, 13).Value = 0.5529953917050692
$ws.Cells.Item(3, 14).Value = 0.5380116959064328
$ws.Cells.Item(3, 15).Value = 0.5314009661835748
$ws.Cells.Item(3, 16).Value = 0.561904761904762
$ws.Cells.Item(3, 17).Value = 0.47
$ws.Cells.Item(3, 18).Value = 0.5327868852459016
$ws.Cells.Item(3, 19).Value = 0.4875621890547264
$ws.Cells.Item(3, 20).Value = 0.5235602094240838
$ws.Cells.Item(3, 21).Value = 0.3842364532019704
$ws.Cells.Item(3, 22).Value = 0.3027027027027027
$ws.Cells.Item(3, 23).Value = 0.3364485981308412
$ws.Cells.Item(3, 24).Value = 0.360655737704918
$ws.Cells.Item(3, 25).Value = 0.4563106796116504
$ws.Cells.Item(3, 26).Value = 0.4079601990049752
$ws.Cells.Item(3, 27).Value = 0.5024630541871921
$ws.Cells.Item(3, 28).Value = 0.4623655913978495
$ws.Cells.Item(3, 29).Value = 0.424778761061947
$ws.Cells.Item(3, 30).Value = 0.4307692307692308
$ws.Cells.Item(3, 31).Value = 0.5043478260869565
# Row 4
$ws.Cells.Item(4, 2).Value = 0.4210526315789473
$ws.Cells.Item(4, 3).Value = 0.5934065934065934
$ws.Cells.Item(4, 4).Value = 0.5305164319248826
$ws.Cells.Item(4, 5).Value = 0.5947368421052631
$ws.Cells.Item(4, 6).Value = 0.6568627450980391
$ws.Cells.Item(4, 7).Value = 0.5050505050505051
$ws.Cells.Item(4, 8).Value = 0.6728110599078341
$ws.Cells.Item(4, 9).Value = 0.5459459459459459
$ws.Cells.Item(4, 10).Value = 0.6683168316831684
$ws.Cells.Item(4, 11).Value = 0.5026737967914439
$ws.Cells.Item(4, 12).Value = 0.5167464114832536
$ws.Cells.Item(4, 13).Value = 0.4908256880733945
$ws.Cells.Item(4, 14).Value = 0.5454545454545454
$ws.Cells.Item(4, 15).Value = 0.5863636363636363
$ws.Cells.Item(4, 16).Value = 0.4687500000000001
$ws.Cells.Item(4, 17).Value = 0.6284153005464481
$ws.Cells.Item(4, 18).Value = 0.4423963133640553
$ws.Cells.Item(4, 19).Value = 0.6382978723404256
$ws.Cells.Item(4, 20).Value = 0.6553398058252428
$ws.Cells.Item(4, 21).Value = 0.5497630331753555
$ws.Cells.Item(4, 22).Value = 0.6380952380952382
# Row 5
$ws.Cells.Item(5, 2).Value = 0.5789473684210527
$ws.Cells.Item(5, 3).Value = 0.4065934065934066
$ws.Cells.Item(5, 4).Value = 0.4694835680751174
$ws.Cells.Item(5, 5).Value = 0.4052631578947369
$ws.Cells.Item(5, 6).Value = 0.3431372549019607
$ws.Cells.Item(5, 7).Value = 0.494949494949495
$ws.Cells.Item(5, 8).Value = 0.3271889400921659
$ws.Cells.Item(5, 9).Value = 0.454054054054054
$ws.Cells.Item(5, 10).Value = 0.3316831683168316
$ws.Cells.Item(5, 11).Value = 0.4973262032085561
$ws.Cells.Item(5, 12).Value = 0.4832535885167464
$ws.Cells.Item(5, 13).Value = 0.5091743119266056
$ws.Cells.Item(5, 14).Value = 0.4545454545454545
$ws.Cells.Item(5, 15).Value = 0.4136363636363636
$ws.Cells.Item(5, 16).Value = 0.53125
$ws.Cells.Item(5, 17).Value = 0.3715846994535519
$ws.Cells.Item(5, 18).Value = 0.5576036866359447
$ws.Cells.Item(5, 19).Value = 0.3617021276595744
$ws.Cells.Item(5, 20).Value = 0.3446601941747572
$ws.Cells.Item(5, 21).Value = 0.4502369668246445
$ws.Cells.Item(5, 22).Value = 0.3619047619047619
# Row 6
$ws.Cells.Item(6, 2).Value = 0.2903225806451613
$ws.Cells.Item(6, 3).Value = 0.297071129707113
$ws.Cells.Item(6, 4).Value = 0.2130857648099027
$ws.Cells.Item(6, 5).Value = 0.2422907488986784
$ws.Cells.Item(6, 6).Value = 0.3220338983050847
$ws.Cells.Item(6, 7).Value = 0.3721973094170403
$ws.Cells.Item(6, 8).Value = 0.3165829145728643
$ws.Cells.Item(6, 9).Value = 0.3409090909090909
$ws.Cells.Item(6, 10).Value = 0.3059360730593607
$ws.Cells.Item(6, 11).Value = 0.3333333333333334
$ws.Cells.Item(6, 12).Value = 0.1380952380952381
$ws.Cells.Item(6, 13).Value = 0.1369047619047619
$ws.Cells.Item(6, 14).Value = 0.1931150293870697
$ws.Cells.Item(6, 15).Value = 0.4044444444444444
$ws.Cells.Item(6, 16).Value = 0.2157598499061914
$ws.Cells.Item(6, 17).Value = 0.1723804925156929
$ws.Cells.Item(6, 18).Value = 0.1146401985111662
$ws.Cells.Item(6, 19).Value = 0.2307692307692308
$ws.Cells.Item(6, 20).Value = 0.1317671092951992
# Row 7
$ws.Cells.Item(7, 2).Value = 0.7096774193548387
$ws.Cells.Item(7, 3).Value = 0.7029288702928871
$ws.Cells.Item(7, 4).Value = 0.7869142351900973
$ws.Cells.Item(7, 5).Value = 0.7577092511013216
$ws.Cells.Item(7, 6).Value = 0.6779661016949152
$ws.Cells.Item(7, 7).Value = 0.6278026905829596
$ws.Cells.Item(7, 8).Value = 0.6834170854271356
$ws.Cells.Item(7, 9).Value = 0.6590909090909091
$ws.Cells.Item(7, 10).Value = 0.6940639269406393
$ws.Cells.Item(7, 11).Value = 0.6666666666666667
$ws.Cells.Item(7, 12).Value = 0.8619047619047618
$ws.Cells.Item(7, 13).Value = 0.8630952380952381
$ws.Cells.Item(7, 14).Value = 0.8068849706129303
$ws.Cells.Item(7, 15).Value = 0.5955555555555556
$ws.Cells.Item(7, 16).Value = 0.7842401500938085
$ws.Cells.Item(7, 17).Value = 0.8276195074843072
$ws.Cells.Item(7, 18).Value = 0.8853598014888338
$ws.Cells.Item(7, 19).Value = 0.7692307692307693
$ws.Cells.Item(7, 20).Value = 0.8682328907048008
